# Applies the edit described by the diff:
#  - Row 7: add F7 = "It gets displayed the Perferences" and G7 = "Pass"
#  - sheetView: selection moves from A15 to F8 (and the topLeftCell scroll anchor
#    is cleared as a natural consequence of re-selecting/saving the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the two new values on row 7 (creates the new shared string automatically)
$ws.Range("F7").Value = "It gets displayed the Perferences"
$ws.Range("G7").Value = "Pass"

# Update the active selection to F8
$ws.Range("F8").Select()
